# TokenIteratorFieldRewriterSplit migration:
# Split the run containing the opening "{" of an m2doc field tag away
# from the rest of the tag text, in the two paragraphs where the tag
# token was still glued together as "{m" / "{m:" in a single w:r.
#
# Word's object model has no "Runs" collection, so runs are (re)created
# implicitly from Range edits. We use Range.InsertXML on a sub-range
# whose End coincides with the paragraph's text end (this keeps the
# edit anchored correctly instead of being appended at the very end of
# the paragraph) to rebuild the remainder of the paragraph as two
# explicit runs: "{" stays untouched before the sub-range, and the new
# XML supplies the "m" (or "m:") run followed by the unchanged tail run.

$d = $word.ActiveDocument

function Split-BraceRun($para, [int]$prefixLen, [string]$firstPiece, [string]$restXml) {
    $full = $para.Range
    $s = $full.Start
    $e = $full.End
    # $e is just past the paragraph mark; the paragraph's own text ends
    # one character earlier.
    $textEnd = $e - 1
    $sub = $d.Range($s + 1, $textEnd)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $firstPiece + '</w:t></w:r>' + $restXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $sub.InsertXML($xml)
}

# --- Edit 1: paragraph styled "Titre1" ("Heading 1"), run text "{m" ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.StartsWith("{m")) {
        $restXml = '<w:r><w:t xml:space="preserve">:v.name}</w:t></w:r>'
        Split-BraceRun $p 1 "m" $restXml
        break
    }
}

# --- Edit 2: paragraph containing "{m:endfor}", run text "{m:" ---
$r = $d.Content
$found = $r.Find.Execute("endfor}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $pos = $r.Start
    for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {
        $p2 = $d.Paragraphs.Item($j)
        if ($pos -ge $p2.Range.Start -and $pos -lt $p2.Range.End) {
            $restXml = '<w:r><w:t xml:space="preserve">endfor}</w:t></w:r>'
            Split-BraceRun $p2 1 "m:" $restXml
            break
        }
    }
}
